$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: Fix the grammar slip "Finally, There is a target field" -> "..., there
# is a target field" and drop the now-unneeded gramStart/gramEnd proofing marks
# that surrounded the capitalized "There".
# ---------------------------------------------------------------------------
$old1 = "distance between points. Finally, There is a target field"
$new1 = "distance between points. Finally, there is a target field"
$null = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# ---------------------------------------------------------------------------
# Hunk 2: relocate the "_GoBack" bookmark from the very end of the document to
# right after "...written using Junit5 and mockito." (end of the decision-tree
# testing paragraph).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $d.Content
$null = $target.Find.Execute("Junit5 and mockito.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
# A zero-length Range confuses Bookmarks.Add, so anchor the bookmark around a
# throwaway character and then shrink the range back down to empty.
$target.InsertAfter("|")
$null = $d.Bookmarks.Add("_GoBack", $target)
$target.Text = ""

# ---------------------------------------------------------------------------
# Hunk 3: merge the two runs "...is working properly. " / "The tests were
# written using Junit5 and " in the RAP testing paragraph into a single run.
# ---------------------------------------------------------------------------
$old3 = "is working properly. The tests were written using Junit5 and "
$null = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2)
